$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format so numeric-looking strings (e.g. "614.68")
# are preserved exactly as text, matching the original inlineStr cell type.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.825.75"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.164.78"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "614.68"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "149.37"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.163.80"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "0.475"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "35.82"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").Value = "3.684.39"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "64.743.79"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "3.160.43"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "6.95"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("D20").Value = "484.06"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "14.78"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "8.07"
$ws.Range("E22").Value = "  +3.74%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "0.717"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "13.78"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "84.28"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("D28").Value = "8.63"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "0.121"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  -6.94%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "26.53"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "1.15"
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").Value = "0.0₃0799"
$ws.Range("E36").Value = "  +7.53%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.30"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "6.06"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "53.31"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").Value = "463.89"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "0.0403"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("D43").Value = "8.43"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "2.857.06"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "2.35"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").Value = "0.271"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  +5.12%  "
$ws.Range("D48").Value = "26.75"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "0.115"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "120.24"
$ws.Range("E51").Value = "  -0.27%  "

# Restore default (unformatted) style on the data range so no stray number format
# is left applied to the cells (keeps cell styling identical to the original file).
$dataRange.ClearFormats()
